$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 629   # was 627
$ws1.Range("F3").Value = 647   # was 638
$ws1.Range("F4").Value = 921   # was 917
$ws1.Range("F6").Value = 819   # was 815
$ws1.Range("F7").Value = 376   # was 374
$ws1.Range("F8").Value = 582   # was 581
$ws1.Range("F9").Value = 120   # was 119
$ws1.Range("F10").Value = 1177   # was 1174
$ws1.Range("F11").Value = 607   # was 602
$ws1.Range("F12").Value = 361   # was 358
$ws1.Range("F13").Value = 486   # was 482
$ws1.Range("F14").Value = 159   # was 158
$ws1.Range("F15").Value = 206   # was 175
$ws1.Range("F16").Value = 322   # was 321
$ws1.Range("F18").Value = 78   # was 77
$ws1.Range("F19").Value = 539   # was 537
$ws1.Range("F20").Value = 53   # was 47
$ws1.Range("F21").Value = 550   # was 548
$ws1.Range("F22").Value = 21   # was 20
$ws1.Range("F23").Value = 586   # was 581

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 77   # was 74
$ws2.Range("F7").Value = 638   # was 637
$ws2.Range("F11").Value = 16   # was 15
$ws2.Range("F13").Value = 52   # was 49
$ws2.Range("G3").Value = "不可售"   # was 78

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 77   # was 74
$ws4.Range("F4").Value = 629   # was 627
$ws4.Range("F7").Value = 647   # was 638
$ws4.Range("F8").Value = 921   # was 917
$ws4.Range("F10").Value = 819   # was 815
$ws4.Range("F11").Value = 376   # was 374
$ws4.Range("F12").Value = 582   # was 581
$ws4.Range("F13").Value = 120   # was 119
$ws4.Range("F14").Value = 1177   # was 1174
$ws4.Range("F15").Value = 607   # was 602
$ws4.Range("F18").Value = 361   # was 358
$ws4.Range("F19").Value = 486   # was 482
$ws4.Range("F20").Value = 638   # was 637
$ws4.Range("F21").Value = 159   # was 158
$ws4.Range("F22").Value = 206   # was 175
$ws4.Range("F24").Value = 322   # was 321
$ws4.Range("F26").Value = 78   # was 77
$ws4.Range("F29").Value = 539   # was 537
$ws4.Range("F30").Value = 16   # was 15
$ws4.Range("F32").Value = 52   # was 49
$ws4.Range("F33").Value = 53   # was 47
$ws4.Range("F34").Value = 550   # was 548
$ws4.Range("F35").Value = 21   # was 20
$ws4.Range("F36").Value = 586   # was 581
$ws4.Range("G5").Value = "不可售"   # was 78
